$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain text so numeric-looking values
# (e.g. "236.02") are not auto-converted to real numbers by Excel,
# matching the original inline-string cell type.
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '25.866.66'
$ws.Range('E2').Value = '  +0.89%  '
$ws.Range('D3').Value = '1.751.69'
$ws.Range('E3').Value = '  +0.74%  '
$ws.Range('E4').Value = '  -0.28%  '
$ws.Range('D5').Value = '236.02'
$ws.Range('E5').Value = '  +0.21%  '
$ws.Range('D6').Value = '0.9985'
$ws.Range('E6').Value = '  -0.27%  '
$ws.Range('D7').Value = '0.5138'
$ws.Range('E7').Value = '  +5.03%  '
$ws.Range('D8').Value = '40.38'
$ws.Range('E8').Value = '  -2.23%  '
$ws.Range('D9').Value = '0.2673'
$ws.Range('E9').Value = '  +4.92%  '
$ws.Range('D10').Value = '0.06178'
$ws.Range('E10').Value = '  +1.80%  '
$ws.Range('D11').Value = '1.753.68'
$ws.Range('E11').Value = '  +0.78%  '
$ws.Range('D12').Value = '0.06943'
$ws.Range('E12').Value = '  +1.59%  '
$ws.Range('D13').Value = '15.47'
$ws.Range('E13').Value = '  +4.91%  '
$ws.Range('D14').Value = '0.6330'
$ws.Range('E14').Value = '  +13.41%  '
$ws.Range('D15').Value = '4.485'
$ws.Range('E15').Value = '  +1.19%  '
$ws.Range('D16').Value = '77.92'
$ws.Range('E16').Value = '  +3.05%  '
$ws.Range('D17').Value = '0.9967'
$ws.Range('E17').Value = '  -0.49%  '
$ws.Range('D18').Value = '0.9989'
$ws.Range('E18').Value = '  -0.19%  '
$ws.Range('D19').Value = '25.895.32'
$ws.Range('E19').Value = '  +0.85%  '
$ws.Range('E20').Value = '  +1.47%  '
$ws.Range('D21').Value = '0.000006661'
$ws.Range('E21').Value = '  +1.99%  '
$ws.Range('D22').Value = '1.997.73'
$ws.Range('E22').Value = '  +1.85%  '
$ws.Range('D23').Value = '4.061'
$ws.Range('E23').Value = '  +0.88%  '
$ws.Range('D24').Value = '8.265'
$ws.Range('E24').Value = '  +4.78%  '
$ws.Range('D25').Value = '5.172'
$ws.Range('E25').Value = '  +3.72%  '
$ws.Range('D26').Value = '136.24'
$ws.Range('E26').Value = '  -0.65%  '
$ws.Range('D27').Value = '1.479'
$ws.Range('E27').Value = '  +0.25%  '
$ws.Range('D28').Value = '15.11'
$ws.Range('E28').Value = '  +3.03%  '
$ws.Range('D29').Value = '1.767'
$ws.Range('E29').Value = '  -2.86%  '
$ws.Range('D30').Value = '102.90'
$ws.Range('E30').Value = '  +2.01%  '
$ws.Range('D32').Value = '3.686'
$ws.Range('E32').Value = '  +0.32%  '
$ws.Range('D33').Value = '3.393'
$ws.Range('E33').Value = '  +0.92%  '
$ws.Range('D34').Value = '0.04388'
$ws.Range('E34').Value = '  -0.25%  '
$ws.Range('D35').Value = '2.641'
$ws.Range('E35').Value = '  +1.19%  '
$ws.Range('D36').Value = '0.9988'
$ws.Range('E36').Value = '  +3.24%  '
$ws.Range('D37').Value = '0.6031'
$ws.Range('E37').Value = '  +2.77%  '
$ws.Range('D38').Value = '2.708'
$ws.Range('E38').Value = '  +2.31%  '
$ws.Range('D39').Value = '0.01560'
$ws.Range('E39').Value = '  +4.23%  '
$ws.Range('D40').Value = '1.935'
$ws.Range('E40').Value = '  +3.77%  '
$ws.Range('D41').Value = '0.9983'
$ws.Range('E41').Value = '  -0.28%  '
$ws.Range('D42').Value = '102.17'
$ws.Range('E42').Value = '  -0.98%  '
$ws.Range('D43').Value = '0.3833'
$ws.Range('E43').Value = '  +3.26%  '
$ws.Range('D44').Value = '0.7484'
$ws.Range('E44').Value = '  +3.66%  '
$ws.Range('D45').Value = '4.897'
$ws.Range('E45').Value = '  -4.43%  '
$ws.Range('D46').Value = '0.05490'
$ws.Range('E46').Value = '  +5.37%  '
$ws.Range('D47').Value = '0.1101'
$ws.Range('E47').Value = '  +2.32%  '
$ws.Range('D48').Value = '5.969'
$ws.Range('E48').Value = '  +3.83%  '
$ws.Range('D49').Value = '30.10'
$ws.Range('E49').Value = '  +0.99%  '
$ws.Range('D50').Value = '52.55'
$ws.Range('E50').Value = '  +1.87%  '
$ws.Range('E51').Value = '  +0.40%  '

# Restore default style on column D so no stray number-format style
# is left attached to the cells (keeps styles.xml usage consistent).
$ws.Range('D2:D51').Style = 'Normal'

